# Weekly update: insert a new record row at row 72 (pushing all
# subsequent rows down by one) and populate it with the new week's
# Apio (Terminal Hortofrutícola Agro Chillán) price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 72:193 down to 73:194, extending the used range to R194.
$ws.Rows.Item(72).Insert()

# Fill in the new row 72 with this week's record.
$ws.Cells.Item(72, 1).Value  = 7
$ws.Cells.Item(72, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value  = "Ñuble"
$ws.Cells.Item(72, 4).Value  = 44662
$ws.Cells.Item(72, 5).Value  = 16
$ws.Cells.Item(72, 6).Value  = 100112017
$ws.Cells.Item(72, 7).Value  = "Apio"
$ws.Cells.Item(72, 8).Value  = "Americana (o)"
$ws.Cells.Item(72, 9).Value  = "Primera"
$ws.Cells.Item(72, 10).Value = 80
$ws.Cells.Item(72, 11).Value = 7500
$ws.Cells.Item(72, 12).Value = 8000
$ws.Cells.Item(72, 13).Value = 7750
$ws.Cells.Item(72, 14).Value = "`$/docena de matas"
$ws.Cells.Item(72, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(72, 16).Value = 1292
$ws.Cells.Item(72, 17).Value = 6
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# Column D carries a custom date/time number format - make sure the
# newly inserted row keeps it (Insert() already copies it from the row
# above, but set it explicitly so the result is format-exact either way).
$ws.Cells.Item(72, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
